$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply an AutoFilter on column E (Team Name, the 5th column of A1:E149)
# keeping only rows whose value is "Team_4" - this hides every other
# data row (the matching rows, plus the header row, stay visible).
$ws.Range("A1:E149").AutoFilter(5, @("Team_4"), 7)

# Matches the author's final cursor position recorded in the saved file.
$ws.Range("K154").Select()
